# Append newly received loss-of-sale records to the Walk-In Report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 22 ----
$ws.Range("A22").Value = 20
$ws.Range("A22").NumberFormat = "0"
$ws.Range("B22").Value = "22-12-2025"
$ws.Range("C22").Value = "athif"
$ws.Range("D22").Value = 9633437635
$ws.Range("D22").NumberFormat = "0"
$ws.Range("E22").Value = "29-12-2025"
$ws.Range("F22").Value = "RASAL"
$ws.Range("G22").Value = "Loss"
$ws.Range("H22").Value = "PRODUCT"
$ws.Range("I22").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J22").Value = "-"
$ws.Range("K22").Value = "BOOTCUT AND DOUBLE BRUST"

# ---- Row 23 ----
$ws.Range("A23").Value = 21
$ws.Range("A23").NumberFormat = "0"
$ws.Range("B23").Value = "23-12-2025"
$ws.Range("C23").Value = "Jemshid"
$ws.Range("D23").Value = 9745172555
$ws.Range("D23").NumberFormat = "0"
$ws.Range("E23").Value = "15-02-2026"
$ws.Range("F23").Value = "RASAL"
$ws.Range("G23").Value = "Loss"
$ws.Range("H23").Value = "ENQUIRY"
$ws.Range("I23").Value = "-"
$ws.Range("J23").Value = "-"

# ---- Row 24 ----
$ws.Range("A24").Value = 22
$ws.Range("A24").NumberFormat = "0"
$ws.Range("B24").Value = "24-12-2025"
$ws.Range("C24").Value = "Adil"
$ws.Range("D24").Value = 7559803880
$ws.Range("D24").NumberFormat = "0"
$ws.Range("E24").Value = "26-01-2026"
$ws.Range("F24").Value = "RASAL"
$ws.Range("G24").Value = "Loss"
$ws.Range("H24").Value = "CUSTOMER INTERNAL ISSUES"
$ws.Range("I24").Value = "FAMILY DISAPPROVEL"
$ws.Range("J24").Value = "-"

Write-Host "Added rows 22-24 to sheet $($ws.Name)"
